$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.127.30"
$ws.Range("E2").Value = "  -3.78%  "

Set-TextValue $ws.Range("D3") "1.849.18"
$ws.Range("E3").Value = "  -2.53%  "

Set-TextValue $ws.Range("D4") "0.9995"
$ws.Range("E4").Value = "  +0.09%  "

Set-TextValue $ws.Range("D5") "0.7075"
$ws.Range("E5").Value = "  -5.54%  "

Set-TextValue $ws.Range("D6") "238.29"
$ws.Range("E6").Value = "  -1.98%  "

Set-TextValue $ws.Range("D7") "0.9997"
$ws.Range("E7").Value = "  +0.06%  "

Set-TextValue $ws.Range("D8") "0.3054"
$ws.Range("E8").Value = "  -3.97%  "

Set-TextValue $ws.Range("D9") "0.07522"
$ws.Range("E9").Value = "  +3.52%  "

Set-TextValue $ws.Range("D10") "23.41"
$ws.Range("E10").Value = "  -6.87%  "

Set-TextValue $ws.Range("D11") "0.08135"
$ws.Range("E11").Value = "  -2.88%  "

Set-TextValue $ws.Range("D12") "1.896.57"
$ws.Range("E12").Value = "  -1.82%  "

Set-TextValue $ws.Range("D13") "0.7255"
$ws.Range("E13").Value = "  -5.13%  "

Set-TextValue $ws.Range("D14") "5.222"
$ws.Range("E14").Value = "  -4.52%  "

Set-TextValue $ws.Range("D15") "89.24"
$ws.Range("E15").Value = "  -4.36%  "

Set-TextValue $ws.Range("D16") "29.263.77"
$ws.Range("E16").Value = "  -3.33%  "

Set-TextValue $ws.Range("D17") "5.798"
$ws.Range("E17").Value = "  -6.60%  "

Set-TextValue $ws.Range("D18") "239.50"
$ws.Range("E18").Value = "  -4.99%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D19") "13.09"
$ws.Range("E19").Value = "  -4.58%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D20") "0.000007682"
$ws.Range("E20").Value = "  -2.67%  "

Set-TextValue $ws.Range("D21") "1.001"
$ws.Range("E21").Value = "  +0.27%  "

Set-TextValue $ws.Range("D22") "2.118.60"
$ws.Range("E22").Value = "  -1.13%  "

Set-TextValue $ws.Range("D23") "0.9998"
$ws.Range("E23").Value = "  +0.14%  "

Set-TextValue $ws.Range("D24") "7.569"
$ws.Range("E24").Value = "  -5.88%  "

Set-TextValue $ws.Range("D25") "0.1469"
$ws.Range("E25").Value = "  -8.49%  "

Set-TextValue $ws.Range("D26") "8.992"
$ws.Range("E26").Value = "  -3.67%  "

Set-TextValue $ws.Range("D27") "161.41"
$ws.Range("E27").Value = "  -1.99%  "

$ws.Range("E28").Value = "  -4.36%  "

Set-TextValue $ws.Range("D29") "1.939"
$ws.Range("E29").Value = "  -6.93%  "

Set-TextValue $ws.Range("D30") "1.385"
$ws.Range("E30").Value = "  -6.29%  "

Set-TextValue $ws.Range("D31") "4.576"
$ws.Range("E31").Value = "  -0.94%  "

$ws.Range("E32").Value = "  -3.04%  "

Set-TextValue $ws.Range("D33") "4.009"
$ws.Range("E33").Value = "  -5.45%  "

Set-TextValue $ws.Range("D34") "0.05178"
$ws.Range("E34").Value = "  -4.87%  "

Set-TextValue $ws.Range("D35") "1.187"
$ws.Range("E35").Value = "  -6.21%  "

Set-TextValue $ws.Range("D36") "1.034"
$ws.Range("E36").Value = "  +3.72%  "

Set-TextValue $ws.Range("D37") "0.7058"
$ws.Range("E37").Value = "  -8.09%  "

Set-TextValue $ws.Range("D38") "2.641"
$ws.Range("E38").Value = "  -2.84%  "

Set-TextValue $ws.Range("D39") "0.01866"
$ws.Range("E39").Value = "  -5.91%  "

Set-TextValue $ws.Range("D40") "2.677"
$ws.Range("E40").Value = "  -3.50%  "

Set-TextValue $ws.Range("D41") "0.9357"
$ws.Range("E41").Value = "  +7.10%  "

Set-TextValue $ws.Range("D42") "6.003"
$ws.Range("E42").Value = "  -1.58%  "

Set-TextValue $ws.Range("D43") "1.079.49"
$ws.Range("E43").Value = "  -2.11%  "

Set-TextValue $ws.Range("D44") "0.4309"
$ws.Range("E44").Value = "  -6.35%  "

Set-TextValue $ws.Range("D45") "70.32"
$ws.Range("E45").Value = "  -4.03%  "

Set-TextValue $ws.Range("D46") "0.9993"
$ws.Range("E46").Value = "  -0.12%  "

Set-TextValue $ws.Range("D47") "102.31"
$ws.Range("E47").Value = "  -2.30%  "

Set-TextValue $ws.Range("D48") "1.753"
$ws.Range("E48").Value = "  -6.79%  "

Set-TextValue $ws.Range("D49") "2.000.22"
$ws.Range("E49").Value = "  -2.68%  "

Set-TextValue $ws.Range("D50") "7.074"
$ws.Range("E50").Value = "  -7.58%  "

Set-TextValue $ws.Range("D51") "9.192"
$ws.Range("E51").Value = "  -4.65%  "

